$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ApiInfoSheet")
$ws2 = $wb.Worksheets.Item("Case2")

# --- Header row text corrections (sheet1 / ApiInfoSheet) ---
# Old headers were the generic/duplicated "ApiId","APIName","Type","Url";
# replace with the descriptive bilingual headers used elsewhere in the workbook.
$ws1.Range("A1").Value = "ApiId(接口编号)"
$ws1.Range("B1").Value = "ApiName(接口名称)"
$ws1.Range("C1").Value = "Type(接口提交方式)"
$ws1.Range("D1").Value = "Url(接口地址)"

# --- Header row text corrections (sheet2 / Case2) ---
$ws2.Range("A1").Value = "CaseId(用例编号)"
$ws2.Range("B1").Value = "ApiId(接口编号)"
$ws2.Range("C1").Value = "Desc(用例描述)"
$ws2.Range("D1").Value = "Params(参数)"

# --- Column width change: sheet1 column C widens from 14 to 23 ---
# (22.15 lands exactly on a stored width of 23 once Excel's internal
# pixel-rounding is applied to the ColumnWidth property.)
$ws1.Columns.Item(3).ColumnWidth = 22.15

# --- Selection / active-tab changes ---
# Update Case2's remembered selection first (without leaving it as the active tab),
# then activate ApiInfoSheet and set its selection - this matches the final state:
# ApiInfoSheet becomes the selected/active tab with C2 selected, while Case2 keeps
# its own last selection at D13 and is no longer the active tab.
[void]$ws2.Select()
[void]$ws2.Range("D13").Select()
[void]$ws1.Select()
[void]$ws1.Range("C2").Select()

Write-Output "done"
